$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: turn the PioneerNoob/Login entry into a new CloneScene/clone entry with ID 3
$ws.Range("A2").Value = "../../NFDataCfg/Ini/NFZoneServer/Scene/CloneScene/"
$ws.Range("B2").Value = "3"
$ws.Range("F2").Value = "clone"

# Rows 3 & 4: rename the SceneName from Stage001 to newscene
$ws.Range("F3").Value = "newscene"
$ws.Range("F4").Value = "newscene"

# Update the active selection to match the saved view state
$ws.Range("H8").Select() | Out-Null
